# faturamento_diario.xlsx — add a missing daily-revenue row for 03/06/2025
# and correct a handful of values pulled from the updated source data.
#
# The sheet is a flat Dia/total_venda/Mes/Ano/Periodo table ordered by
# Periodo then Dia. A new day (Dia=3, Mes=6) belongs right after the
# existing Dia=2/Mes=6 row, so every row from the old row 3 down to the
# old last row (62) shifts down by one. A few total_venda values in the
# 05/2025 block also needed correcting once the new data came in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing 02/06/2025 total (was 3672.65).
$ws.Range("B2").Value = 34060.21

# Make room for the new 03/06/2025 row by shifting rows 3:62 down to 4:63.
$ws.Rows.Item(3).Insert()

# New row: Dia=3, total_venda=4281.25, Mes=6, Ano=2025, Periodo=06/2025
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 4281.25
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 2025
$ws.Range("E3").Value = "06/2025"

# A few total_venda corrections inside the (now shifted) 05/2025 block.
$ws.Range("B20").Value = 27080.43
$ws.Range("B23").Value = 28806.18
$ws.Range("B24").Value = 18949.46
